# Applies the re-ordering of rows 2-11 (by columns A,B,D,E,F,G,H,Q,R)
# as described by the diff. The other columns (C,I,P,S,T,U,V,W,Y,Z,AA,AB,
# AD,AE,AG,AT,AW,AX,AY) are identical across all these rows, so only the
# changed columns need to be rewritten with their new (reordered) values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @{ Row = 2;  A = 88042571; B = 89356; D = "LC"; E = 5447; F = "Vedticka";       G = "Fuscoporia viticola";        H = "(Schwein.) Murrill";         Q = 401467.8944965442; R = 6660836.946491199 },
    @{ Row = 3;  A = 88042580; B = 94121; D = "NT"; E = 53;   F = "Vedtrappmossa";  G = "Crossocalyx hellerianus";    H = "(Nees ex Lindenb.) Meyl.";   Q = 401426.0083459655; R = 6660767.15908254 },
    @{ Row = 4;  A = 88042577; B = 77506; D = "NT"; E = 6425; F = "Garnlav";        G = "Alectoria sarmentosa";       H = "(Ach.) Ach.";                Q = 401472.0542245907; R = 6660936.206280074 },
    @{ Row = 5;  A = 88042570; B = 94121; D = "NT"; E = 53;   F = "Vedtrappmossa";  G = "Crossocalyx hellerianus";    H = "(Nees ex Lindenb.) Meyl.";   Q = 401456.0475465701; R = 6660786.82821779 },
    @{ Row = 6;  A = 88042584; B = 94121; D = "NT"; E = 53;   F = "Vedtrappmossa";  G = "Crossocalyx hellerianus";    H = "(Nees ex Lindenb.) Meyl.";   Q = 401460.8001688122; R = 6660796.188836097 },
    @{ Row = 7;  A = 88042575; B = 94121; D = "NT"; E = 53;   F = "Vedtrappmossa";  G = "Crossocalyx hellerianus";    H = "(Nees ex Lindenb.) Meyl.";   Q = 401585.1169067804; R = 6660693.984798764 },
    @{ Row = 8;  A = 88042573; B = 94121; D = "NT"; E = 53;   F = "Vedtrappmossa";  G = "Crossocalyx hellerianus";    H = "(Nees ex Lindenb.) Meyl.";   Q = 401628.0818258527; R = 6660804.195007879 },
    @{ Row = 9;  A = 88042578; B = 94121; D = "NT"; E = 53;   F = "Vedtrappmossa";  G = "Crossocalyx hellerianus";    H = "(Nees ex Lindenb.) Meyl.";   Q = 401610.2002120143; R = 6660790.191765637 },
    @{ Row = 10; A = 88042587; B = 94121; D = "NT"; E = 53;   F = "Vedtrappmossa";  G = "Crossocalyx hellerianus";    H = "(Nees ex Lindenb.) Meyl.";   Q = 402072.8906266145; R = 6660546.13208588 },
    @{ Row = 11; A = 88042582; B = 94121; D = "NT"; E = 53;   F = "Vedtrappmossa";  G = "Crossocalyx hellerianus";    H = "(Nees ex Lindenb.) Meyl.";   Q = 402104.0025583604; R = 6660549.798474666 }
)

foreach ($item in $data) {
    $r = $item.Row
    $ws.Cells.Item($r, 1).Value = $item.A   # A
    $ws.Cells.Item($r, 2).Value = $item.B   # B
    $ws.Cells.Item($r, 4).Value = $item.D   # D
    $ws.Cells.Item($r, 5).Value = $item.E   # E
    $ws.Cells.Item($r, 6).Value = $item.F   # F
    $ws.Cells.Item($r, 7).Value = $item.G   # G
    $ws.Cells.Item($r, 8).Value = $item.H   # H
    $ws.Cells.Item($r, 17).Value = $item.Q  # Q
    $ws.Cells.Item($r, 18).Value = $item.R  # R
}
